$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename existing "New CRM bottle" entries (F28:F29) to note the 02/07 bottle-open date ---
$ws.Range("F28").Value = "New CRM bottle (opened 02/07)"
$ws.Range("F29").Value = "New CRM bottle (opened 02/07)"

# --- Row 30: new CRM-bracketed sample measured 02/14 (date serial 43510) ---
# Copy the date format (s="1", numFmtId 14 = m/d/yyyy) from an existing date cell so the
# new cell picks up the identical style index instead of minting a new one.
$ws.Range("A29").Copy()
$ws.Range("A30").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A30").Value = 43510
$ws.Range("B30").Value = 2219.88
$ws.Range("F30").Value = "New CRM bottle (opened 02/14)"

# --- Row 31: keep it empty, but carry the date style forward like the rest of column A ---
$ws.Range("A29").Copy()
$ws.Range("A31").PasteSpecial(-4122)  # xlPasteFormats

# --- Column F width bump (cosmetic, to fit the longer label text) ---
$ws.Columns("F").ColumnWidth = 29.666666666666668

# --- View state: scrolled down a few rows further, selection moved to F37 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
$ws.Range("F37").Select() | Out-Null

$excel.CutCopyMode = $false
